# Numeracion.xlsx - "move LOS calculation after the filling of blank dates"
# On Hoja1, rows 8-10 (Age/Birthdate/Gender block) are rotated by one position:
#   old row 10 (Gender)    -> row 8
#   old row 8  (Age)       -> row 9
#   old row 9  (Birthdate) -> row 10
# Capture the old values first, then write them back in the new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$row8 = @($ws.Range("B8").Value(), $ws.Range("C8").Value(), $ws.Range("D8").Value())
$row9 = @($ws.Range("B9").Value(), $ws.Range("C9").Value(), $ws.Range("D9").Value())
$row10 = @($ws.Range("B10").Value(), $ws.Range("C10").Value(), $ws.Range("D10").Value())

# New row 8 = old row 10
$ws.Range("B8").Value = $row10[0]
$ws.Range("C8").Value = $row10[1]
$ws.Range("D8").Value = $row10[2]

# New row 9 = old row 8
$ws.Range("B9").Value = $row8[0]
$ws.Range("C9").Value = $row8[1]
$ws.Range("D9").Value = $row8[2]

# New row 10 = old row 9
$ws.Range("B10").Value = $row9[0]
$ws.Range("C10").Value = $row9[1]
$ws.Range("D10").Value = $row9[2]

# Update the sheet selection to match the authored state.
$ws.Range("A9").Select() | Out-Null
